$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.509.04'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '2.037.74'
$ws.Range('E3').Value = '  +2.97%  '
$ws.Range('E4').Value = '  +0.04%  '
$c = $ws.Range('D5')
$c.Value = '''254.89'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +3.88%  '
$c = $ws.Range('D6')
$c.Value = '''0.622'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -1.42%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range('D7')
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('B8').Value = 'Solana'
$ws.Range('C8').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D8')
$c.Value = '''57.44'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -6.27%  '
$c = $ws.Range('D9')
$c.Value = '''0.387'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +1.29%  '
$c = $ws.Range('D10')
$c.Value = '''57.07'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -0.59%  '
$c = $ws.Range('D11')
$c.Value = '''0.0801'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('E12').Value = '  -0.87%  '
$c = $ws.Range('D13')
$c.Value = '''14.83'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +2.12%  '
$ws.Range('D14').Value = '2.344.61'
$ws.Range('E14').Value = '  +3.38%  '
$c = $ws.Range('D15')
$c.Value = '''0.823'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -2.94%  '
$c = $ws.Range('D16')
$c.Value = '''21.47'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -2.86%  '
$c = $ws.Range('D17')
$c.Value = '''5.38'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -0.81%  '
$ws.Range('D18').Value = '2.042.42'
$ws.Range('E18').Value = '  +3.39%  '
$ws.Range('D19').Value = '37.431.39'
$ws.Range('E19').Value = '  +1.95%  '
$c = $ws.Range('D20')
$c.Value = '''70.01'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').Value = '0.0₃0857'
$ws.Range('E21').Value = '  -0.49%  '
$c = $ws.Range('D22')
$c.Value = '''5.25'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +2.08%  '
$c = $ws.Range('D23')
$c.Value = '''229.10'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.57%  '
$c = $ws.Range('D24')
$c.Value = '''2.66'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +7.04%  '
$c = $ws.Range('D25')
$c.Value = '''0.999'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('E27').Value = '  -3.79%  '
$c = $ws.Range('D28')
$c.Value = '''9.16'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -1.24%  '
$c = $ws.Range('D29')
$c.Value = '''163.27'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +0.10%  '
$c = $ws.Range('D30')
$c.Value = '''19.90'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +2.13%  '
$c = $ws.Range('D31')
$c.Value = '''1.35'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('E32').Value = '  -0.68%  '
$c = $ws.Range('D33')
$c.Value = '''4.77'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -1.45%  '
$c = $ws.Range('D34')
$c.Value = '''0.0665'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +6.97%  '
$c = $ws.Range('D35')
$c.Value = '''4.50'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -0.85%  '
$c = $ws.Range('D36')
$c.Value = '''2.48'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +8.91%  '
$c = $ws.Range('D37')
$c.Value = '''3.49'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +3.82%  '
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('E39').Value = '  +2.12%  '
$c = $ws.Range('D40')
$c.Value = '''5.39'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('E41').Value = '  +4.17%  '
$c = $ws.Range('D42')
$c.Value = '''0.0971'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -0.35%  '
$c = $ws.Range('D43')
$c.Value = '''0.0218'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +2.93%  '
$c = $ws.Range('D44')
$c.Value = '''1.19'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +1.30%  '
$c = $ws.Range('D45')
$c.Value = '''16.34'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +1.40%  '
$ws.Range('D46').Value = '1.402.46'
$ws.Range('E46').Value = '  +2.25%  '
$c = $ws.Range('D47')
$c.Value = '''91.42'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +1.91%  '
$c = $ws.Range('D48')
$c.Value = '''1.05'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.19%  '
$c = $ws.Range('D49')
$c.Value = '''7.44'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +2.78%  '
$c = $ws.Range('D50')
$c.Value = '''2.88'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('E51').Value = '  +7.03%  '
